$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 94. This shifts existing rows 94:212 down to 95:213,
# carrying their values/formatting with them (matching the diff where every
# row from 95..213 now holds what used to be in the row above it).
$ws.Rows("94:94").Insert()

# Populate the newly inserted row 94 with the new data record.
$ws.Range("A94").Value = 8
$ws.Range("B94").Value = "Terminal La Palmera de La Serena"
$ws.Range("C94").Value = "Coquimbo"
$ws.Range("D94").Value = 45174
$ws.Range("D94").NumberFormat = $ws.Range("D95").NumberFormat
$ws.Range("E94").Value = 4
$ws.Range("F94").Value = 100112052
$ws.Range("G94").Value = "Albahaca"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = 3250
$ws.Range("N94").Value = "`$/paquete"
$ws.Range("O94").Value = "Región de Arica y Parinacota"
$ws.Range("P94").Value = 3250
$ws.Range("Q94").Value = 1
$ws.Range("R94").Value = "Hortaliza"
